{"js": "// Replace the date line and each \"NNN\u00d7N=\" equation in the practice sheet\n// with the values from the new day's worksheet. Every old string occurs\n// exactly once in the document, so a scoped, case-sensitive search +\n// in-place replace for each pair is safe and preserves run formatting.\nconst replacements = [\n  [\"2024-10-19 Saturday\", \"2024-10-20 Sunday\"],\n  [\"906\u00d79=\", \"466\u00d74=\"],\n  [\"831\u00d79=\", \"801\u00d75=\"],\n  [\"320\u00d74=\", \"683\u00d74=\"],\n  [\"356\u00d75=\", \"740\u00d73=\"],\n  [\"895\u00d76=\", \"523\u00d74=\"],\n  [\"506\u00d75=\", \"674\u00d75=\"],\n  [\"863\u00d75=\", \"389\u00d77=\"],\n  [\"661\u00d73=\", \"696\u00d76=\"],\n  [\"429\u00d79=\", \"197\u00d79=\"],\n  [\"361\u00d73=\", \"269\u00d79=\"],\n  [\"363\u00d75=\", \"389\u00d72=\"],\n  [\"724\u00d76=\", \"484\u00d73=\"],\n  [\"635\u00d73=\", \"105\u00d72=\"],\n  [\"831\u00d74=\", \"517\u00d73=\"],\n  [\"681\u00d72=\", \"327\u00d77=\"],\n  [\"896\u00d79=\", \"481\u00d77=\"],\n  [\"909\u00d77=\", \"422\u00d74=\"],\n  [\"885\u00d73=\", \"331\u00d78=\"],\n  [\"772\u00d72=\", \"741\u00d75=\"],\n  [\"977\u00d72=\", \"701\u00d74=\"],\n  [\"155\u00d78=\", \"220\u00d75=\"],\n  [\"746\u00d74=\", \"981\u00d77=\"],\n  [\"259\u00d74=\", \"244\u00d76=\"],\n  [\"280\u00d77=\", \"299\u00d72=\"],\n  [\"335\u00d73=\", \"817\u00d79=\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each \"NNN\u00d7N=\" equation in the practice sheet\n# with the values from the new day's worksheet. Every old string occurs\n# exactly once in the document, so a document-wide Find/Replace (wrap\n# continue, replace-all) for each pair is safe and preserves run formatting.\n$replacements = @(\n  @{ old = \"2024-10-19 Saturday\"; new = \"2024-10-20 Sunday\" },\n  @{ old = \"906\u00d79=\"; new = \"466\u00d74=\" },\n  @{ old = \"831\u00d79=\"; new = \"801\u00d75=\" },\n  @{ old = \"320\u00d74=\"; new = \"683\u00d74=\" },\n  @{ old = \"356\u00d75=\"; new = \"740\u00d73=\" },\n  @{ old = \"895\u00d76=\"; new = \"523\u00d74=\" },\n  @{ old = \"506\u00d75=\"; new = \"674\u00d75=\" },\n  @{ old = \"863\u00d75=\"; new = \"389\u00d77=\" },\n  @{ old = \"661\u00d73=\"; new = \"696\u00d76=\" },\n  @{ old = \"429\u00d79=\"; new = \"197\u00d79=\" },\n  @{ old = \"361\u00d73=\"; new = \"269\u00d79=\" },\n  @{ old = \"363\u00d75=\"; new = \"389\u00d72=\" },\n  @{ old = \"724\u00d76=\"; new = \"484\u00d73=\" },\n  @{ old = \"635\u00d73=\"; new = \"105\u00d72=\" },\n  @{ old = \"831\u00d74=\"; new = \"517\u00d73=\" },\n  @{ old = \"681\u00d72=\"; new = \"327\u00d77=\" },\n  @{ old = \"896\u00d79=\"; new = \"481\u00d77=\" },\n  @{ old = \"909\u00d77=\"; new = \"422\u00d74=\" },\n  @{ old = \"885\u00d73=\"; new = \"331\u00d78=\" },\n  @{ old = \"772\u00d72=\"; new = \"741\u00d75=\" },\n  @{ old = \"977\u00d72=\"; new = \"701\u00d74=\" },\n  @{ old = \"155\u00d78=\"; new = \"220\u00d75=\" },\n  @{ old = \"746\u00d74=\"; new = \"981\u00d77=\" },\n  @{ old = \"259\u00d74=\"; new = \"244\u00d76=\" },\n  @{ old = \"280\u00d77=\"; new = \"299\u00d72=\" },\n  @{ old = \"335\u00d73=\"; new = \"817\u00d79=\" }\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $pair.old\n  $find.Replacement.Text = $pair.new\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n  # wdFindContinue=1, wdReplaceAll=2\n  $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n}\n"}
